$wb = $excel.ActiveWorkbook

# Sheet ALC row 7
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 10000
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -10224

# Sheet ALC row 9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1819
$ws.Range("I9").Value = 2217.5
$ws.Range("J9").Value = 225
$ws.Range("K9").Value = 2217.5
$ws.Range("L9").Value = 225
$ws.Range("M9").Value = -2048.5
$ws.Range("N9").Value = -563

# Sheet ALC row 14
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 10000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -10382

# Sheet ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 976.1875
$ws.Range("I129").Value = 215.2
$ws.Range("J129").Value = 1040.678
$ws.Range("K129").Value = 645.5999999999999
$ws.Range("L129").Value = 3122.034000000001
$ws.Range("M129").Value = 4354.4
$ws.Range("N129").Value = -13122.034

# Sheet ARM row 41
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 3500
$ws.Range("I41").Value = 3500
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 3500
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -3086
$ws.Range("N41").ClearContents()

# Sheet ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 437502.16
$ws.Range("I61").Value = 386703.72
$ws.Range("J61").Value = 503540.1
$ws.Range("K61").Value = 386703.72
$ws.Range("L61").Value = 503540.1
$ws.Range("M61").Value = -386491.72
$ws.Range("N61").Value = -503964.1

# Sheet ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 213239.31
$ws.Range("I74").Value = 271298.1
$ws.Range("J74").Value = 70027.60000000001
$ws.Range("K74").Value = 271298.1
$ws.Range("L74").Value = 70027.60000000001
$ws.Range("M74").Value = -270424.1
$ws.Range("N74").Value = -71775.60000000001

# Sheet ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 213239.31
$ws.Range("I77").Value = 271298.1
$ws.Range("J77").Value = 70027.60000000001
$ws.Range("K77").Value = 1356490.5
$ws.Range("L77").Value = 350138
$ws.Range("M77").Value = -1352122.5
$ws.Range("N77").Value = -358874

# Sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 22545.42
$ws.Range("I132").Value = 28628.842
$ws.Range("J132").Value = 3281.25
$ws.Range("K132").Value = 85886.526
$ws.Range("L132").Value = 9843.75
$ws.Range("M132").Value = -83356.526
$ws.Range("N132").Value = -14903.75

# Sheet ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 437502.16
$ws.Range("I136").Value = 386703.72
$ws.Range("J136").Value = 503540.1
$ws.Range("K136").Value = 1160111.16
$ws.Range("L136").Value = 1510620.3
$ws.Range("M136").Value = -1157561.16
$ws.Range("N136").Value = -1515720.3

# Sheet BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3802.1667
$ws.Range("I99").Value = 4027.4375
$ws.Range("K99").Value = 4027.4375
$ws.Range("M99").Value = -2529.4375

# Sheet BSM row 126
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H126").Value = 27500
$ws.Range("J126").Value = 27500
$ws.Range("L126").Value = 27500
$ws.Range("N126").Value = -37380

# Sheet BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1869.0975
$ws.Range("I134").Value = 1029.75
$ws.Range("J134").Value = 3676.923
$ws.Range("K134").Value = 3089.25
$ws.Range("L134").Value = 11030.769
$ws.Range("M134").Value = -554.25
$ws.Range("N134").Value = -16100.769

# Sheet CRP row 3
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 18248.25
$ws.Range("I3").Value = 125
$ws.Range("J3").Value = 36371.5
$ws.Range("K3").Value = 125
$ws.Range("L3").Value = 36371.5
$ws.Range("M3").Value = -12
$ws.Range("N3").Value = -36597.5

# Sheet CRP row 21
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 16000
$ws.Range("J21").Value = 16000
$ws.Range("L21").Value = 16000
$ws.Range("N21").Value = -16470

# Sheet CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3202.361
$ws.Range("I31").Value = 2253.2693
$ws.Range("J31").Value = 5670
$ws.Range("K31").Value = 2253.2693
$ws.Range("L31").Value = 5670
$ws.Range("M31").Value = -1958.2693
$ws.Range("N31").Value = -6260

# Sheet CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3202.361
$ws.Range("I34").Value = 2253.2693
$ws.Range("J34").Value = 5670
$ws.Range("K34").Value = 2253.2693
$ws.Range("L34").Value = 5670
$ws.Range("M34").Value = -2051.2693
$ws.Range("N34").Value = -6074

# Sheet CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 821.3333
$ws.Range("I105").Value = 808.8889
$ws.Range("K105").Value = 808.8889
$ws.Range("M105").Value = 938.1111

# Sheet CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1434.8889
$ws.Range("I122").Value = 1233.3334
$ws.Range("J122").Value = 1838
$ws.Range("K122").Value = 3700.0002
$ws.Range("L122").Value = 5514
$ws.Range("M122").Value = -1250.0002
$ws.Range("N122").Value = -10414

# Sheet CUL row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 769313.9399999999
$ws.Range("I2").Value = 2000036
$ws.Range("J2").Value = 112.625
$ws.Range("K2").Value = 12000216
$ws.Range("L2").Value = 675.75
$ws.Range("M2").Value = -12000103
$ws.Range("N2").Value = -901.75

# Sheet CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1371.9
$ws.Range("I5").Value = 399.7143
$ws.Range("J5").Value = 1895.3846
$ws.Range("K5").Value = 1199.1429
$ws.Range("L5").Value = 5686.1538
$ws.Range("M5").Value = -1087.1429
$ws.Range("N5").Value = -5910.1538

# Sheet CUL row 20
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 600
$ws.Range("J20").Value = 600
$ws.Range("L20").Value = 1800
$ws.Range("N20").Value = -2254

# Sheet CUL row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 1000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 6000
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -6566

# Sheet CUL row 64
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 1811.9445
$ws.Range("I64").Value = 851.5
$ws.Range("J64").Value = 3012.5
$ws.Range("K64").Value = 2554.5
$ws.Range("L64").Value = 9037.5
$ws.Range("M64").Value = -2284.5
$ws.Range("N64").Value = -9577.5

# Sheet CUL row 67
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H67").Value = 1811.9445
$ws.Range("I67").Value = 851.5
$ws.Range("J67").Value = 3012.5
$ws.Range("K67").Value = 2554.5
$ws.Range("L67").Value = 9037.5
$ws.Range("M67").Value = -1618.5
$ws.Range("N67").Value = -10909.5

# Sheet CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1371.9
$ws.Range("I135").Value = 399.7143
$ws.Range("J135").Value = 1895.3846
$ws.Range("K135").Value = 3597.4287
$ws.Range("L135").Value = 17058.4614
$ws.Range("M135").Value = -1062.4287
$ws.Range("N135").Value = -22128.4614

# Sheet GSM row 140
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 29950
$ws.Range("J140").Value = 29950
$ws.Range("L140").Value = 29950
$ws.Range("N140").Value = -40310

# Sheet LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1568.6666
$ws.Range("I16").Value = 1641.2858
$ws.Range("J16").Value = 1162
$ws.Range("K16").Value = 1641.2858
$ws.Range("L16").Value = 1162
$ws.Range("M16").Value = -1471.2858
$ws.Range("N16").Value = -1502

# Sheet LTW row 24
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 3502
$ws.Range("J24").Value = 3500
$ws.Range("L24").Value = 3500
$ws.Range("N24").Value = -4186

# Sheet LTW row 35
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 3635.4412
$ws.Range("I35").Value = 1718
$ws.Range("J35").Value = 3966.0344
$ws.Range("K35").Value = 1718
$ws.Range("L35").Value = 3966.0344
$ws.Range("M35").Value = -1382
$ws.Range("N35").Value = -4638.0344

# Sheet LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1512.3334
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# Sheet LTW row 94
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 26553.334
$ws.Range("J94").Value = 26553.334
$ws.Range("L94").Value = 26553.334
$ws.Range("N94").Value = -27905.334

# Sheet LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 76928856
$ws.Range("I100").Value = 8842.857
$ws.Range("K100").Value = 8842.857
$ws.Range("M100").Value = -8301.857

# Sheet LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1512.3334
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# Sheet WVR row 3
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 201530
$ws.Range("I3").Value = 667000
$ws.Range("J3").Value = 2042.8572
$ws.Range("K3").Value = 667000
$ws.Range("L3").Value = 2042.8572
$ws.Range("M3").Value = -666886
$ws.Range("N3").Value = -2270.8572

# Sheet WVR row 19
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 20377.75
$ws.Range("I19").Value = 3505
$ws.Range("J19").Value = 26002
$ws.Range("K19").Value = 3505
$ws.Range("L19").Value = 26002
$ws.Range("M19").Value = -3331
$ws.Range("N19").Value = -26350

# Sheet WVR row 31
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 3933.5
$ws.Range("J31").Value = 4000
$ws.Range("L31").Value = 4000
$ws.Range("N31").Value = -4696

# Sheet WVR row 70
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 18333.334
$ws.Range("I70").Value = 5000
$ws.Range("J70").Value = 25000
$ws.Range("K70").Value = 5000
$ws.Range("L70").Value = 25000
$ws.Range("M70").Value = -4685
$ws.Range("N70").Value = -25630

# Sheet WVR row 73
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 18333.334
$ws.Range("I73").Value = 5000
$ws.Range("J73").Value = 25000
$ws.Range("K73").Value = 5000
$ws.Range("L73").Value = 25000
$ws.Range("M73").Value = -3908
$ws.Range("N73").Value = -27184

# Sheet WVR row 80
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 36150.5
$ws.Range("J80").Value = 36150.5
$ws.Range("L80").Value = 36150.5
$ws.Range("N80").Value = -38146.5

# Sheet WVR row 83
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H83").Value = 36150.5
$ws.Range("J83").Value = 36150.5
$ws.Range("L83").Value = 108451.5
$ws.Range("N83").Value = -118435.5

# Sheet WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1051.5
$ws.Range("I126").Value = 786.95
$ws.Range("J126").Value = 1933.3334
$ws.Range("K126").Value = 2360.85
$ws.Range("L126").Value = 5800.0002
$ws.Range("M126").Value = 109.1499999999996
$ws.Range("N126").Value = -10740.0002
